# Auto-generated Excel COM-interop script to apply Titan_Profits.xlsx diff
# Updates currentAveragePrice/LevePrice/LeveProfit columns (H:N) for specific rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 361478.6
$ws.Range("I64").Value = 628393.8
$ws.Range("J64").Value = 5591.6665
$ws.Range("K64").Value = 628393.8
$ws.Range("L64").Value = 5591.6665
$ws.Range("M64").Value = -628145.8
$ws.Range("N64").Value = -6087.6665

$ws.Range("H67").Value = 361478.6
$ws.Range("I67").Value = 628393.8
$ws.Range("J67").Value = 5591.6665
$ws.Range("K67").Value = 628393.8
$ws.Range("L67").Value = 5591.6665
$ws.Range("M67").Value = -627535.8
$ws.Range("N67").Value = -7307.6665

$ws.Range("H74").Value = 4015.3845
$ws.Range("I74").Value = 3662.5
$ws.Range("J74").Value = 4580
$ws.Range("K74").Value = 3662.5
$ws.Range("L74").Value = 4580
$ws.Range("M74").Value = -2726.5
$ws.Range("N74").Value = -6452

$ws.Range("H77").Value = 4015.3845
$ws.Range("I77").Value = 3662.5
$ws.Range("J77").Value = 4580
$ws.Range("K77").Value = 18312.5
$ws.Range("L77").Value = 22900
$ws.Range("M77").Value = -13632.5
$ws.Range("N77").Value = -32260

$ws.Range("H96").Value = 472.4762
$ws.Range("I96").Value = 396.2
$ws.Range("J96").Value = 541.8182
$ws.Range("K96").Value = 1188.6
$ws.Range("L96").Value = 1625.4546
$ws.Range("M96").Value = 184.4000000000001
$ws.Range("N96").Value = -4371.4546

$ws.Range("H100").Value = 16668371
$ws.Range("I100").Value = 33335616
$ws.Range("J100").Value = 1126
$ws.Range("K100").Value = 33335616
$ws.Range("L100").Value = 1126
$ws.Range("M100").Value = -33335075

$ws.Range("H138").Value = 4169491.8
$ws.Range("I138").Value = 1159684.8
$ws.Range("J138").Value = 6176030
$ws.Range("K138").Value = 3479054.4
$ws.Range("L138").Value = 18528090
$ws.Range("M138").Value = -3473914.4
$ws.Range("N138").Value = -18538370

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13751.298
$ws.Range("I32").Value = 1796.1025
$ws.Range("J32").Value = 169168.83
$ws.Range("K32").Value = 1796.1025
$ws.Range("L32").Value = 169168.83
$ws.Range("M32").Value = -1509.1025
$ws.Range("N32").Value = -169742.83

$ws.Range("H97").Value = 55572490
$ws.Range("I97").Value = 66686908
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 66686908
$ws.Range("L97").Value = 400
$ws.Range("M97").Value = -66686412
$ws.Range("N97").Value = -1392

$ws.Range("H102").Value = 1960
$ws.Range("I102").Value = 1950
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1950
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -328
$ws.Range("N102").Value = -5244

$ws.Range("H122").Value = 2539.389
$ws.Range("I122").Value = 2332.9167
$ws.Range("J122").Value = 2952.3333
$ws.Range("K122").Value = 6998.750100000001
$ws.Range("L122").Value = 8856.999899999999
$ws.Range("M122").Value = -4548.750100000001
$ws.Range("N122").Value = -13756.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5664.0835
$ws.Range("I86").Value = 1571.875
$ws.Range("J86").Value = 13848.5
$ws.Range("K86").Value = 1571.875
$ws.Range("L86").Value = 13848.5
$ws.Range("M86").Value = -448.875
$ws.Range("N86").Value = -16094.5

$ws.Range("H89").Value = 5664.0835
$ws.Range("I89").Value = 1571.875
$ws.Range("J89").Value = 13848.5
$ws.Range("K89").Value = 7859.375
$ws.Range("L89").Value = 69242.5
$ws.Range("M89").Value = -2243.375
$ws.Range("N89").Value = -80474.5

$ws.Range("H94").Value = 1332.2354
$ws.Range("I94").Value = 1300.9
$ws.Range("J94").Value = 1377
$ws.Range("K94").Value = 1300.9
$ws.Range("L94").Value = 1377
$ws.Range("M94").Value = -849.9000000000001

$ws.Range("H99").Value = 2677.625
$ws.Range("I99").Value = 2627.5
$ws.Range("J99").Value = 2727.75
$ws.Range("K99").Value = 2627.5
$ws.Range("L99").Value = 2727.75
$ws.Range("M99").Value = -1129.5
$ws.Range("N99").Value = -5723.75

$ws.Range("H105").Value = 3090.6
$ws.Range("I105").Value = 2959.476
$ws.Range("J105").Value = 3287.2856
$ws.Range("K105").Value = 2959.476
$ws.Range("L105").Value = 3287.2856
$ws.Range("M105").Value = -1212.476
$ws.Range("N105").Value = -6781.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2147.1538
$ws.Range("I31").Value = 1161.4667
$ws.Range("J31").Value = 3491.2727
$ws.Range("K31").Value = 1161.4667
$ws.Range("L31").Value = 3491.2727
$ws.Range("M31").Value = -866.4666999999999
$ws.Range("N31").Value = -4081.2727

$ws.Range("H34").Value = 2147.1538
$ws.Range("I34").Value = 1161.4667
$ws.Range("J34").Value = 3491.2727
$ws.Range("K34").Value = 1161.4667
$ws.Range("L34").Value = 3491.2727
$ws.Range("M34").Value = -959.4666999999999
$ws.Range("N34").Value = -3895.2727

$ws.Range("H62").Value = 26722.223
$ws.Range("I62").Value = 36083.332
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 36083.332
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -35459.332

$ws.Range("H65").Value = 26722.223
$ws.Range("I65").Value = 36083.332
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 180416.66
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -177296.66

$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 5000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -4251

$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 15000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -11256

$ws.Range("H105").Value = 1011
$ws.Range("I105").Value = 866.61536
$ws.Range("J105").Value = 1949.5
$ws.Range("K105").Value = 866.61536
$ws.Range("L105").Value = 1949.5
$ws.Range("M105").Value = 880.38464

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 949.875
$ws.Range("I138").Value = 949.875
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2849.625
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 2290.375
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 6318.3335
$ws.Range("I141").Value = 10303.333
$ws.Range("J141").Value = 2333.3333
$ws.Range("K141").Value = 30909.999
$ws.Range("L141").Value = 6999.999899999999
$ws.Range("M141").Value = -25729.999
$ws.Range("N141").Value = -17359.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 730.625
$ws.Range("I97").Value = 705.7692
$ws.Range("J97").Value = 760
$ws.Range("K97").Value = 705.7692
$ws.Range("L97").Value = 760
$ws.Range("M97").Value = -209.7692
$ws.Range("N97").Value = -1752

$ws.Range("H122").Value = 794918
$ws.Range("I122").Value = 855927.0600000001
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 2567781.18
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -2565331.18
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 686.25
$ws.Range("I16").Value = 798
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 798
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -628
$ws.Range("N16").Value = -840

$ws.Range("H55").Value = 672.75
$ws.Range("I55").Value = 649.5
$ws.Range("J55").Value = 680.5
$ws.Range("K55").Value = 649.5
$ws.Range("L55").Value = 680.5
$ws.Range("M55").Value = -476.5
$ws.Range("N55").Value = -1026.5

$ws.Range("H68").Value = 2159.6667
$ws.Range("I68").Value = 2066.25
$ws.Range("J68").Value = 2284.2222
$ws.Range("K68").Value = 2066.25
$ws.Range("L68").Value = 2284.2222
$ws.Range("M68").Value = -1317.25
$ws.Range("N68").Value = -3782.2222

$ws.Range("H71").Value = 2159.6667
$ws.Range("I71").Value = 2066.25
$ws.Range("J71").Value = 2284.2222
$ws.Range("K71").Value = 10331.25
$ws.Range("L71").Value = 11421.111
$ws.Range("M71").Value = -6587.25
$ws.Range("N71").Value = -18909.111

$ws.Range("H82").Value = 1235.3
$ws.Range("I82").Value = 1100
$ws.Range("J82").Value = 1269.125
$ws.Range("K82").Value = 1100
$ws.Range("L82").Value = 1269.125
$ws.Range("M82").Value = -739
$ws.Range("N82").Value = -1991.125

$ws.Range("H85").Value = 1235.3
$ws.Range("I85").Value = 1100
$ws.Range("J85").Value = 1269.125
$ws.Range("K85").Value = 1100
$ws.Range("L85").Value = 1269.125
$ws.Range("M85").Value = 148
$ws.Range("N85").Value = -3765.125

$ws.Range("H93").Value = 1416.4736
$ws.Range("I93").Value = 1175.875
$ws.Range("J93").Value = 1591.4546
$ws.Range("K93").Value = 1175.875
$ws.Range("L93").Value = 1591.4546
$ws.Range("M93").Value = 72.125
$ws.Range("N93").Value = -4087.4546

$ws.Range("H100").Value = 2512.739
$ws.Range("I100").Value = 1832.3334
$ws.Range("J100").Value = 2752.8823
$ws.Range("K100").Value = 1832.3334
$ws.Range("L100").Value = 2752.8823
$ws.Range("M100").Value = -1291.3334
$ws.Range("N100").Value = -3834.8823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7706586.5
$ws.Range("I62").Value = 14304986
$ws.Range("J62").Value = 8454
$ws.Range("K62").Value = 14304986
$ws.Range("L62").Value = 8454
$ws.Range("M62").Value = -14304362
$ws.Range("N62").Value = -9702

$ws.Range("H65").Value = 7706586.5
$ws.Range("I65").Value = 14304986
$ws.Range("J65").Value = 8454
$ws.Range("K65").Value = 71524930
$ws.Range("L65").Value = 42270
$ws.Range("M65").Value = -71521810
$ws.Range("N65").Value = -48510

$ws.Range("H81").Value = 627644.25
$ws.Range("I81").Value = 1112728.5
$ws.Range("J81").Value = 3964.4285
$ws.Range("K81").Value = 2225457
$ws.Range("L81").Value = 7928.857
$ws.Range("M81").Value = -2224396
$ws.Range("N81").Value = -10050.857

$ws.Range("H84").Value = 627644.25
$ws.Range("I84").Value = 1112728.5
$ws.Range("J84").Value = 3964.4285
$ws.Range("K84").Value = 11127285
$ws.Range("L84").Value = 39644.285
$ws.Range("M84").Value = -11121981
$ws.Range("N84").Value = -50252.285

$ws.Range("H122").Value = 79231.38
$ws.Range("I122").Value = 93273.45
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 279820.35
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -277370.35

$ws.Range("H136").Value = 9037001
$ws.Range("I136").Value = 10132174
$ws.Range("J136").Value = 1822.5
$ws.Range("K136").Value = 30396522
$ws.Range("L136").Value = 5467.5
$ws.Range("M136").Value = -30393972
